$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in SNR header and values in column P (rows 5-8)
$ws.Range("P5").Value = "SNR"
$ws.Range("P6").Value = 5
$ws.Range("P7").Value = 11
$ws.Range("P8").Value = 11

# Widen column O (the 15th column) to match new content
$ws.Columns.Item(15).ColumnWidth = 22.67

# Move the active selection to P9 (reflecting the new view state)
$ws.Range("P9").Select()
